# Update cryptocurrency price and 1h-volume-change figures
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '26.677.24'
$ws.Cells.Item(2, 5).Value = '  -0.11%  '

$ws.Cells.Item(3, 4).Value = '1.595.53'
$ws.Cells.Item(3, 5).Value = '  -0.25%  '

$ws.Cells.Item(4, 5).Value = '  +0.38%  '

$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = '211.24'
$cell.Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  -0.01%  '

$ws.Cells.Item(6, 5).Value = '  -0.36%  '

$ws.Cells.Item(7, 5).Value = '  +0.36%  '

$ws.Cells.Item(8, 5).Value = '  -0.11%  '

$ws.Cells.Item(9, 5).Value = '  +0.22%  '

$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = '19.41'
$cell.Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  -1.04%  '

$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0841'
$cell.Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  -0.26%  '

$ws.Cells.Item(12, 4).Value = '1.821.26'
$ws.Cells.Item(12, 5).Value = '  -0.13%  '

$ws.Cells.Item(13, 4).Value = '1.592.94'
$ws.Cells.Item(13, 5).Value = '  -1.07%  '

$ws.Cells.Item(14, 5).Value = '  +0.04%  '

$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.522'
$cell.Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  +0.24%  '

$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = '64.97'
$cell.Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  +0.05%  '

$ws.Cells.Item(17, 4).Value = '26.644.60'
$ws.Cells.Item(17, 5).Value = '  -0.16%  '

$ws.Cells.Item(18, 4).Value = '0.0₃0749'
$ws.Cells.Item(18, 5).Value = '  +2.80%  '

$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.01'
$cell.Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  +0.36%  '

$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = '208.93'
$cell.Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  +0.20%  '

$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = '6.97'
$cell.Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  +2.59%  '

$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = '4.27'
$cell.Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  +0.46%  '

$ws.Cells.Item(23, 5).Value = '  -0.45%  '

$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = '142.83'
$cell.Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  -2.06%  '

$ws.Cells.Item(26, 5).Value = '  +0.33%  '

$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = '7.11'
$cell.Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  -1.66%  '

$ws.Cells.Item(28, 5).Value = '  -1.19%  '

$ws.Cells.Item(29, 5).Value = '  +0.21%  '

$ws.Cells.Item(30, 5).Value = '  +1.60%  '

$ws.Cells.Item(31, 5).Value = '  -0.12%  '

$ws.Cells.Item(32, 5).Value = '  +0.26%  '

$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.94'
$cell.Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  +0.25%  '

$ws.Cells.Item(34, 4).Value = '1.284.98'
$ws.Cells.Item(34, 5).Value = '  -0.29%  '

$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.616'
$cell.Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  -7.39%  '

$ws.Cells.Item(36, 5).Value = '  -0.23%  '

$ws.Cells.Item(37, 5).Value = '  -0.61%  '

$ws.Cells.Item(38, 5).Value = '  -0.34%  '

$ws.Cells.Item(39, 5).Value = '  +16.98%  '

$ws.Cells.Item(40, 5).Value = '  -2.14%  '

$ws.Cells.Item(41, 5).Value = '  +0.04%  '

$ws.Cells.Item(42, 5).Value = '  -0.58%  '

$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.781'
$cell.Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  -0.68%  '

$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = '62.91'
$cell.Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  -1.05%  '

$ws.Cells.Item(45, 4).Value = '1.732.84'
$ws.Cells.Item(45, 5).Value = '  -0.21%  '

$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = '90.71'
$cell.Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  +0.73%  '

$ws.Cells.Item(47, 5).Value = '  -3.13%  '

$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.100'
$cell.Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  -0.80%  '

$ws.Cells.Item(49, 5).Value = '  +0.68%  '

$ws.Cells.Item(50, 5).Value = '  +0.31%  '

$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = '7.28'
$cell.Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  -2.75%  '
